# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# 1. Insert a new "Player Info" sheet before "ODI Batting".
# 2. On "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE and replace the
#    full scorecard URL with just the numeric match code; drop the stray
#    empty INNING_NUMBER cells for the "did not bat" rows.
# 3. Append a new "ODI Batting Extra" sheet after "ODI Batting" holding the
#    newly scraped per-match stats.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "ODI Batting" edits (existing sheet)
# ---------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")

# Header rename
$batting.Cells.Item(1, 4).Value = "MATCH_CODE"

# Replace MATCH_CARD_LINK urls with the bare MatchCode number
for ($r = 2; $r -le 56; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $text = $cell.Text
    if ($text -match 'MatchCode=(\d+)') {
        $cell.NumberFormat = "@"
        $cell.Value = $matches[1]
    }
}

# Drop the leftover empty INNING_NUMBER cells on "did not bat" rows
$batting.Cells.Item(5, 2).ClearContents()
$batting.Cells.Item(36, 2).ClearContents()
$batting.Cells.Item(43, 2).ClearContents()

# ---------------------------------------------------------------------
# Step 2: new "Player Info" sheet, inserted before "ODI Batting"
# ---------------------------------------------------------------------
$battingRef = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingRef)
$playerInfo.Name = "Player Info"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    $cell.Value = $piHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$playerInfo.Cells.Item(2, 1).NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value = "4316"
$playerInfo.Cells.Item(2, 2).Value = "Dickwella Patabendige Dilantha Niroshan Dickwella"
$playerInfo.Cells.Item(2, 3).Value = "Left Handed"
$playerInfo.Cells.Item(2, 4).Value = "Does Not Bowl | Unknown"

# ---------------------------------------------------------------------
# Step 3: new "ODI Batting Extra" sheet, inserted after "ODI Batting"
# ---------------------------------------------------------------------
$battingRef2 = $wb.Worksheets.Item("ODI Batting")
$extra = $wb.Worksheets.Add($null, $battingRef2)
$extra.Name = "ODI Batting Extra"

$exHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $exHeaders.Length; $c++) {
    $cell = $extra.Cells.Item(1, $c)
    $cell.Value = $exHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$exRows = @(
    @("4124", "4", "4", "0", "19.00%", "NO"),
    @("4182", "1", "0", "0", "1.04%", "NO"),
    @("4183", $null, $null, $null, $null, "NO"),
    @("4186", "2", "2", "0", "3.51%", "NO"),
    @("4187", "1", "6", "0", "11.11%", "NO"),
    @("4188", $null, $null, $null, $null, "NO"),
    @("4209", $null, $null, $null, $null, "NO"),
    @("4210", "1", "2", "0", "6.43%", "NO"),
    @("4211", "1", "8", "0", "24.00%", "NO"),
    @("4212", "1", "5", "0", "19.05%", "NO"),
    @("4215", "1", "12", "0", "25.96%", "YES"),
    @("4231", "1", "8", "3", "23.31%", "NO"),
    @("4232", "1", "1", "0", "3.02%", "NO"),
    @("4233", $null, $null, $null, $null, "NO"),
    @("4261", "1", "2", "0", "3.46%", "NO"),
    @("4264", "1", "0", "0", "4.35%", "NO"),
    @("4269", "2", "0", "0", "1.65%", "NO"),
    @("4465", "6", "0", "0", "2.45%", "NO"),
    @("4600", $null, $null, $null, $null, "NO"),
    @("4601", "1", "0", "0", "0.39%", "NO")
)

$r = 2
foreach ($row in $exRows) {
    $extra.Cells.Item($r, 1).NumberFormat = "@"
    $extra.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne $null) {
        $extra.Cells.Item($r, 2).Value = [double]$row[1]
    }
    if ($row[2] -ne $null) {
        $extra.Cells.Item($r, 3).NumberFormat = "@"
        $extra.Cells.Item($r, 3).Value = $row[2]
    }
    if ($row[3] -ne $null) {
        $extra.Cells.Item($r, 4).NumberFormat = "@"
        $extra.Cells.Item($r, 4).Value = $row[3]
    }
    if ($row[4] -ne $null) {
        $extra.Cells.Item($r, 5).NumberFormat = "@"
        $extra.Cells.Item($r, 5).Value = $row[4]
    }
    $extra.Cells.Item($r, 6).Value = $row[5]
    $r++
}
